# Updated cryptos list on Sat Jun 15 22:11:26 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "66.095.91";   E = "  +0.19%  " },
    @{ Row = 3;  D = "3.559.66";    E = "  +2.26%  " },
    @{ Row = 4;  D = $null;         E = "  +0.08%  " },
    @{ Row = 5;  D = "606.22";      E = "  +0.65%  " },
    @{ Row = 6;  D = "145.00";      E = "  +1.63%  " },
    @{ Row = 7;  D = "3.557.51";    E = "  +2.26%  " },
    @{ Row = 9;  D = $null;         E = "  +3.39%  " },
    @{ Row = 10; D = $null;         E = "  +1.12%  " },
    @{ Row = 11; D = "7.92";        E = "  -2.83%  " },
    @{ Row = 12; D = "0.412";       E = "  +0.02%  " },
    @{ Row = 13; D = "4.166.40";    E = "  +2.21%  " },
    @{ Row = 14; D = "0.0000207";   E = "  +2.23%  " },
    @{ Row = 15; D = "30.01";       E = "  -0.73%  " },
    @{ Row = 16; D = "3.559.04";    E = "  +1.96%  " },
    @{ Row = 17; D = "66.231.24";   E = "  +0.21%  " },
    @{ Row = 18; D = $null;         E = "  -0.95%  " },
    @{ Row = 19; D = "11.47";       E = "  +10.20%  " },
    @{ Row = 20; D = "6.21";        E = "  +0.59%  " },
    @{ Row = 21; D = $null;         E = "  +1.09%  " },
    @{ Row = 22; D = "429.61";      E = "  +2.33%  " },
    @{ Row = 23; D = "0.614";       E = "  +4.76%  " },
    @{ Row = 24; D = "79.12";       E = "  +1.92%  " },
    @{ Row = 25; D = "3.704.20";    E = "  +2.03%  " },
    @{ Row = 27; D = "0.0000118";   E = "  +1.93%  " },
    @{ Row = 28; D = "2.50";        E = "  +1.81%  " },
    @{ Row = 29; D = "7.95";        E = "  +0.18%  " },
    @{ Row = 30; D = "9.12";        E = "  -3.30%  " },
    @{ Row = 31; D = $null;         E = "  -0.01%  " },
    @{ Row = 32; D = "25.57";       E = "  +1.93%  " },
    @{ Row = 33; D = $null;         E = "  -1.26%  " },
    @{ Row = 34; D = "3.559.25";    E = "  +2.24%  " },
    @{ Row = 35; D = "0.152";       E = "  -6.16%  " },
    @{ Row = 36; D = $null;         E = "  +0.06%  " },
    @{ Row = 37; D = $null;         E = "  +2.19%  " },
    @{ Row = 38; D = "7.87";        E = "  +3.06%  " },
    @{ Row = 39; D = "5.60";        E = "  +0.89%  " },
    @{ Row = 40; D = $null;         E = "  -0.02%  " },
    @{ Row = 41; D = "173.89";      E = "  +1.95%  " },
    @{ Row = 42; D = "0.0849";      E = "  -1.99%  " },
    @{ Row = 43; D = "5.21";        E = "  +2.09%  " },
    @{ Row = 44; D = "0.896";       E = "  +0.42%  " },
    @{ Row = 45; D = $null;         E = "  +1.78%  " },
    @{ Row = 46; D = "46.09";       E = "  +0.96%  " },
    @{ Row = 47; D = "25.85";       E = "  -0.87%  " },
    @{ Row = 48; D = $null;         E = "  -0.26%  " },
    @{ Row = 49; D = $null;         E = "  +1.20%  " },
    @{ Row = 50; D = "23.63";       E = "  +10.11%  " },
    @{ Row = 51; D = "7.13";        E = "  +0.23%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
